$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'71.326.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.841.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'706.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'173.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.92%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.840.04"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.77%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.47%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.26%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'36.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.58%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.489.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.71%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.919.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'71.391.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.08%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.55%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.50%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'498.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.733"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.33%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'85.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'10.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.90%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.994.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.76%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.81%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.42%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'29.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.179"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'9.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.53%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.805.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.12%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.19%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.103"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.21%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +4.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'3.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.11%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.000318"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.51%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'164.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'431.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.33%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'49.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.67%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.19%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.03%  "
$ws.Range("E51").Style = "Normal"
